$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the data rows for 2008 and 2009 (rows 2 and 3), which shifts
# the remaining rows (2010-2015) up to occupy rows 2-7.
$ws.Range("A2:G3").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
